# GRFS.xlsx update: shift the Model sheet's P&L block down by three rows
# (making room above it), and leave the Model sheet as the active/selected
# sheet with the new top row selected - matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model")

# Insert three blank rows above the existing data block (old row 3 becomes
# row 6, etc.) - formulas/relative refs shift automatically.
$ws.Rows("3:5").Insert()

# The freshly inserted rows picked up the surrounding column formatting;
# clear them so they go back to being genuinely blank (no cell records).
$ws.Range("A3:J5").Clear()

# Make "Model" the active sheet/tab, and select the new top data row
# (matches the saved selection in the workbook).
$ws.Activate()
$ws.Range("A3:XFD3").Select()
